$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "dnasr281@gmail.com, "
$prefixLen = $prefix.Length

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.StartsWith($prefix)) {
        $rest = $val.Substring($prefixLen)
        $cell.Value = $rest + ", dnasr281@gmail.com"
    }
}
